$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark that currently sits right
#    after "agents" (before "/agentsService") in the service-location
#    paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Delete the whole paragraph that starts with "צריך להדגיש כי ...",
#    which merges it away and leaves the following paragraph
#    ("כמו כן ישנה אפשרות ...") directly after the previous one.
#    (Range.Paragraphs.First.Range can report just the matched text,
#    so resolve the real paragraph index first and grab it from
#    Document.Paragraphs, whose .Range spans the whole paragraph
#    including its end-of-paragraph mark.)
# ------------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute("צריך להדגיש כי", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $paraIndex = $find.Paragraphs.First.Index
    $para = $d.Paragraphs.Item($paraIndex)
    $para.Range.Delete()
}

# ------------------------------------------------------------------
# 3. Re-insert the "_GoBack" bookmark inside the "כמו כן ..." run,
#    splitting it into "כמ" and "ו כן ישנה אפשרות שמבנה ה ".
# ------------------------------------------------------------------
$find2 = $d.Content
$found2 = $find2.Find.Execute("כמו כן ישנה", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $splitPoint = $find2.Start + 2
    $bmRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
